$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "shape 1"
$ws.Range("B2").Value = ""
$ws.Range("D2").Value = "2021091611DocumentBeek1.xlsx"
$ws.Range("E2").Value = "C:\excel-files\2021091611DocumentBeek1.xlsx"
